$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- TextBox 33: move up slightly (y offset 16226329 -> 16022645 EMU) ---
$tb33 = $s.Shapes.Item("TextBox 33")
$tb33.Top = 1261.6255905511812

# --- TextBox 38: replace the lorem-ipsum paragraph with real text, and shrink height ---
$tb38 = $s.Shapes.Item("TextBox 38")
$tf38 = $tb38.TextFrame
$tf38.DeleteText()
$tf38.TextRange.Text = "We then created front end web pages based on our Figma prototype using HMTL, JavaScript and CSS. The accessibility features were added using JavaScript functions tied to buttons, normally using the onclick() method to trigger a change to the browsers local CSS file. "
$tf38.TextRange.LanguageID = "en-GB"
$tb38.Height = 207.1112
